$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceCell($cellRef, $value) {
    $ws.Range($cellRef).Value = "'" + $value
    $ws.Range($cellRef).Style = "Normal"
}

function Set-TextCell($cellRef, $value) {
    $ws.Range($cellRef).Value = $value
}

# Row 2 - Bitcoin
Set-PriceCell "D2" "62.889.20"
Set-TextCell "E2" "  +4.74%  "

# Row 3 - Ethereum
Set-PriceCell "D3" "3.361.19"
Set-TextCell "E3" "  +5.14%  "

# Row 4 - TetherUSD
Set-TextCell "E4" "  +0.03%  "

# Row 5 - BNB
Set-PriceCell "D5" "560.24"
Set-TextCell "E5" "  +4.04%  "

# Row 6 - Solana
Set-PriceCell "D6" "153.42"
Set-TextCell "E6" "  +5.80%  "

# Row 7 - USDC
Set-TextCell "E7" "  -0.05%  "

# Row 8 - XRP
Set-TextCell "E8" "  +1.02%  "

# Row 9 - Toncoin
Set-PriceCell "D9" "7.54"
Set-TextCell "E9" "  +2.67%  "

# Row 10 - Dogecoin
Set-TextCell "E10" "  +4.29%  "

# Row 11 - Cardano
Set-PriceCell "D11" "0.439"
Set-TextCell "E11" "  +1.88%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-PriceCell "D12" "3.940.27"
Set-TextCell "E12" "  +5.19%  "

# Row 13 - TRON
Set-TextCell "E13" "  +0.35%  "

# Row 14/15 - ShibaInu/Avalanche swap content
Set-TextCell "B14" "Avalanche"
Set-TextCell "C14" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-PriceCell "D14" "27.14"
Set-TextCell "E14" "  +4.19%  "

Set-TextCell "B15" "ShibaInu"
Set-TextCell "C15" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-PriceCell "D15" "0.0000182"
Set-TextCell "E15" "  +3.73%  "

# Row 16 - WrappedBTC
Set-PriceCell "D16" "62.997.63"
Set-TextCell "E16" "  +4.87%  "

# Row 17 - WrappedEther
Set-PriceCell "D17" "3.357.36"
Set-TextCell "E17" "  +4.78%  "

# Row 18 - Polkadot
Set-PriceCell "D18" "6.51"
Set-TextCell "E18" "  +4.39%  "

# Row 19 - Chainlink
Set-PriceCell "D19" "13.84"
Set-TextCell "E19" "  +5.68%  "

# Row 20 - Uniswap
Set-PriceCell "D20" "8.45"
Set-TextCell "E20" "  +1.27%  "

# Row 21 - BitcoinCash
Set-PriceCell "D21" "389.88"
Set-TextCell "E21" "  +1.72%  "

# Row 22 - Polygon
Set-PriceCell "D22" "0.543"
Set-TextCell "E22" "  +2.46%  "

# Row 23 - Dai
Set-TextCell "E23" "  +0.23%  "

# Row 24 - Litecoin
Set-TextCell "E24" "  +0.34%  "

# Row 25 - Kaspa
Set-TextCell "E25" "  +5.13%  "

# Row 26 - InternetComputer(DFINITY)
Set-PriceCell "D26" "8.87"
Set-TextCell "E26" "  +0.11%  "

# Row 27 - PEPE
Set-PriceCell "D27" "0.0₃0972"
Set-TextCell "E27" "  +7.48%  "

# Row 29 - RenderToken
Set-PriceCell "D29" "6.69"
Set-TextCell "E29" "  +8.24%  "

# Row 30 - PancakeSwap
Set-PriceCell "D30" "1.99"
Set-TextCell "E30" "  +4.67%  "

# Row 31 - NEARProtocol
Set-PriceCell "D31" "5.68"
Set-TextCell "E31" "  +5.07%  "

# Row 32 - EthereumClassic
Set-PriceCell "D32" "23.09"
Set-TextCell "E32" "  +3.14%  "

# Row 33 - Fetch.AI
Set-PriceCell "D33" "1.31"
Set-TextCell "E33" "  +6.90%  "

# Row 34 - Aptos
Set-TextCell "E34" "  +1.81%  "

# Row 35 - ImmutableX
Set-PriceCell "D35" "1.48"
Set-TextCell "E35" "  +9.25%  "

# Row 36 - Monero
Set-PriceCell "D36" "160.71"
Set-TextCell "E36" "  +2.80%  "

# Row 37 - Stacks
Set-PriceCell "D37" "1.89"
Set-TextCell "E37" "  +12.17%  "

# Row 38 - EnergySwap
Set-PriceCell "D38" "27.02"
Set-TextCell "E38" "  +4.73%  "

# Row 39 - Hedera
Set-TextCell "E39" "  +4.57%  "

# Row 40 - Maker
Set-PriceCell "D40" "2.830.96"
Set-TextCell "E40" "  +1.77%  "

# Row 41 - VeChain
Set-PriceCell "D41" "0.0310"
Set-TextCell "E41" "  +8.33%  "

# Row 42 - Filecoin
Set-PriceCell "D42" "4.34"
Set-TextCell "E42" "  +2.16%  "

# Row 43/44 - Mantle/OKB swap content
Set-TextCell "B43" "OKB"
Set-TextCell "C43" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-PriceCell "D43" "40.79"
Set-TextCell "E43" "  +2.54%  "

Set-TextCell "B44" "Mantle"
Set-TextCell "C44" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-PriceCell "D44" "0.749"
Set-TextCell "E44" "  +2.95%  "

# Row 45 - ONDO
Set-PriceCell "D45" "1.05"
Set-TextCell "E45" "  +5.20%  "

# Row 46 - InjectiveProtocol
Set-PriceCell "D46" "22.23"
Set-TextCell "E46" "  +8.20%  "

# Row 47 - RenzoRestakedETH
Set-PriceCell "D47" "3.405.77"
Set-TextCell "E47" "  +5.24%  "

# Row 48 - Stellar
Set-TextCell "E48" "  +2.56%  "

# Row 49 - Cosmos
Set-TextCell "E49" "  +2.41%  "

# Row 50 - SuiNetwork
Set-PriceCell "D50" "0.809"
Set-TextCell "E50" "  +0.71%  "

# Row 51 - Bittensor
Set-PriceCell "D51" "282.28"
Set-TextCell "E51" "  +6.32%  "
